$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.57779315203136
$ws.Range("C2").Value = 7.153282327285898
$ws.Range("D2").Value = 15.01852548308429
$ws.Range("E2").Value = 16.43938439622382
$ws.Range("G2").Value = 3.678728604813291
$ws.Range("I2").Value = 24.68942862365345
$ws.Range("J2").Value = 9.389190241959341
$ws.Range("K2").Value = 12.10415100787166
$ws.Range("N2").Value = 19.89365978858281
$ws.Range("O2").Value = 27.91304022677859
$ws.Range("B3").Value = 11.33747745152123
$ws.Range("C3").Value = 6.991311992468935
$ws.Range("D3").Value = 14.95722939437158
$ws.Range("E3").Value = 16.37832576323501
$ws.Range("G3").Value = 3.680809086434632
$ws.Range("I3").Value = 24.76923467071082
$ws.Range("J3").Value = 9.396768550694265
$ws.Range("K3").Value = 11.9476152764383
$ws.Range("N3").Value = 19.95638184321639
$ws.Range("O3").Value = 27.97775780621648
$ws.Range("B4").Value = 11.18965504868252
$ws.Range("C4").Value = 6.891459503336968
$ws.Range("D4").Value = 14.9228615470779
$ws.Range("E4").Value = 16.34444537351779
$ws.Range("G4").Value = 3.682154500080106
$ws.Range("I4").Value = 24.82261680847136
$ws.Range("J4").Value = 9.402827071894752
$ws.Range("K4").Value = 11.85278262364351
$ws.Range("N4").Value = 19.99665944274593
$ws.Range("O4").Value = 28.02302374536799
$ws.Range("B5").Value = 11.12943650631922
$ws.Range("C5").Value = 6.850728562821762
$ws.Range("D5").Value = 14.9096889192517
$ws.Range("E5").Value = 16.33155681974135
$ws.Range("G5").Value = 3.682719918455314
$ws.Range("I5").Value = 24.84547078565421
$ws.Range("J5").Value = 9.405649603452352
$ws.Range("K5").Value = 11.81450423532104
$ws.Range("N5").Value = 20.01351832565
$ws.Range("O5").Value = 28.042857771399
$ws.Range("B6").Value = 11.11944108813534
$ws.Range("C6").Value = 6.843964605607996
$ws.Range("D6").Value = 14.9075521910359
$ws.Range("E6").Value = 16.32947240652293
$ws.Range("G6").Value = 3.682814843262408
$ws.Range("I6").Value = 24.84933208987644
$ws.Range("J6").Value = 9.406139646626542
$ws.Range("K6").Value = 11.80817157729564
$ws.Range("N6").Value = 20.01634467336626
$ws.Range("O6").Value = 28.04623493712883
$ws.Range("B7").Value = 11.1888427209865
$ws.Range("C7").Value = 6.89091027277922
$ws.Range("D7").Value = 14.92268051206498
$ws.Range("E7").Value = 16.34426782439266
$ws.Range("G7").Value = 3.682162055991837
$ws.Range("I7").Value = 24.82292057177871
$ws.Range("J7").Value = 9.40286370543784
$ws.Range("K7").Value = 11.85226484471505
$ws.Range("N7").Value = 19.99688500220175
$ws.Range("O7").Value = 28.02328561841656
$ws.Range("B8").Value = 11.49503968413103
$ws.Range("C8").Value = 7.097554298995037
$ws.Range("D8").Value = 14.99671903506671
$ws.Range("E8").Value = 16.41758863799294
$ws.Range("G8").Value = 3.679431876734613
$ws.Range("I8").Value = 24.71603572348615
$ws.Range("J8").Value = 9.391511711586102
$ws.Range("K8").Value = 12.0499385295201
$ws.Range("N8").Value = 19.91492064033809
$ws.Range("O8").Value = 27.93420581195162
$ws.Range("B9").Value = 12.08964566269534
$ws.Range("C9").Value = 7.497027546906139
$ws.Range("D9").Value = 15.16732062061111
$ws.Range("E9").Value = 16.58951274092629
$ws.Range("G9").Value = 3.674614966408677
$ws.Range("I9").Value = 24.54125258366309
$ws.Range("J9").Value = 9.380389340527689
$ws.Range("K9").Value = 12.44576324200518
$ws.Range("N9").Value = 19.76813809661353
$ws.Range("O9").Value = 27.80349879447746
$ws.Range("B10").Value = 12.51817657042624
$ws.Range("C10").Value = 7.78376031433471
$ws.Range("D10").Value = 15.30739179142947
$ws.Range("E10").Value = 16.73224179037228
$ws.Range("G10").Value = 3.671399825848347
$ws.Range("I10").Value = 24.43414276371128
$ws.Range("J10").Value = 9.3789882864462
$ws.Range("K10").Value = 12.73891912491138
$ws.Range("N10").Value = 19.66871080711513
$ws.Range("O10").Value = 27.73441584603839
$ws.Range("B11").Value = 12.71038579108592
$ws.Range("C11").Value = 7.912106192313038
$ws.Range("D11").Value = 15.37413067701567
$ws.Range("E11").Value = 16.80055919200647
$ws.Range("G11").Value = 3.670006749443488
$ws.Range("I11").Value = 24.39005786070342
$ws.Range("J11").Value = 9.379815125561159
$ws.Range("K11").Value = 12.87223802954111
$ws.Range("N11").Value = 19.62528649794849
$ws.Range("O11").Value = 27.70886459467104
$ws.Range("B12").Value = 12.78270572064729
$ws.Range("C12").Value = 7.960358951437732
$ws.Range("D12").Value = 15.39982059298566
$ws.Range("E12").Value = 16.82690011096179
$ws.Range("G12").Value = 3.669489166169251
$ws.Range("I12").Value = 24.37403261986602
$ws.Range("J12").Value = 9.380338096043266
$ws.Range("K12").Value = 12.92267156888187
$ws.Range("N12").Value = 19.609101072157
$ws.Range("O12").Value = 27.70003542826225
$ws.Range("B13").Value = 12.767152143151
$ws.Range("C13").Value = 7.949983114917305
$ws.Range("D13").Value = 15.39426952507841
$ws.Range("E13").Value = 16.82120646828836
$ws.Range("G13").Value = 3.669600195502775
$ws.Range("I13").Value = 24.37745417488983
$ws.Range("J13").Value = 9.380216143099592
$ws.Range("K13").Value = 12.91181280808276
$ws.Range("N13").Value = 19.61257541954625
$ws.Range("O13").Value = 27.70189927245898
$ws.Range("B14").Value = 12.71634532525298
$ws.Range("C14").Value = 7.916083236419416
$ws.Range("D14").Value = 15.37623593697422
$ws.Range("E14").Value = 16.80271694265236
$ws.Range("G14").Value = 3.669963968541105
$ws.Range("I14").Value = 24.38872604291443
$ws.Range("J14").Value = 9.379853948910267
$ws.Range("K14").Value = 12.87638851618242
$ws.Range("N14").Value = 19.62394974118814
$ws.Range("O14").Value = 27.70812123897902
$ws.Range("B15").Value = 12.68516197696608
$ws.Range("C15").Value = 7.89527174744117
$ws.Range("D15").Value = 15.3652436805584
$ws.Range("E15").Value = 16.79145234889781
$ws.Range("G15").Value = 3.670188083608868
$ws.Range("I15").Value = 24.39571752449794
$ws.Range("J15").Value = 9.379659403168365
$ws.Range("K15").Value = 12.85468207668353
$ws.Range("N15").Value = 19.63095045876264
$ws.Range("O15").Value = 27.71204266423234
$ws.Range("B16").Value = 12.50555395645498
$ws.Range("C16").Value = 7.775326326465285
$ws.Range("D16").Value = 15.30308958737408
$ws.Range("E16").Value = 16.72784387344404
$ws.Range("G16").Value = 3.671492261001336
$ws.Range("I16").Value = 24.43711731502187
$ws.Range("J16").Value = 9.37896365742276
$ws.Range("K16").Value = 12.73020169114224
$ws.Range("N16").Value = 19.6715849166391
$ws.Range("O16").Value = 27.73620405286695
$ws.Range("B17").Value = 12.39461741549521
$ws.Range("C17").Value = 7.701172753461023
$ws.Range("D17").Value = 15.26572172952874
$ws.Range("E17").Value = 16.68967853785173
$ws.Range("G17").Value = 3.672310097970959
$ws.Range("I17").Value = 24.46370424191061
$ws.Range("J17").Value = 9.378911396314068
$ws.Range("K17").Value = 12.65379402880807
$ws.Range("N17").Value = 19.69697443246382
$ws.Range("O17").Value = 27.75253230613349
$ws.Range("B18").Value = 12.3305565987086
$ws.Range("C18").Value = 7.65832751230931
$ws.Range("D18").Value = 15.24451402618721
$ws.Range("E18").Value = 16.66804691665684
$ws.Range("G18").Value = 3.672787041502346
$ws.Range("I18").Value = 24.47943297705609
$ws.Range("J18").Value = 9.379019200012168
$ws.Range("K18").Value = 12.60984535021059
$ws.Range("N18").Value = 19.71174784112754
$ws.Range("O18").Value = 27.76247677430086
$ws.Range("B19").Value = 12.30882553606244
$ws.Range("C19").Value = 7.643789051410711
$ws.Range("D19").Value = 15.23738296523033
$ws.Range("E19").Value = 16.66077829094339
$ws.Range("G19").Value = 3.672949652101489
$ws.Range("I19").Value = 24.48483340197684
$ws.Range("J19").Value = 9.37907939787692
$ws.Range("K19").Value = 12.59496632308305
$ws.Range("N19").Value = 19.71677910420028
$ws.Range("O19").Value = 27.76593870257824
$ws.Range("B20").Value = 12.40645353107962
$ws.Range("C20").Value = 7.709086975062027
$ws.Range("D20").Value = 15.26967019431396
$ws.Range("E20").Value = 16.69370829571996
$ws.Range("G20").Value = 3.672222360778834
$ws.Range("I20").Value = 24.46082881644764
$ws.Range("J20").Value = 9.37890269596001
$ws.Range("K20").Value = 12.66192819508801
$ws.Range("N20").Value = 19.69425408658191
$ws.Range("O20").Value = 27.7507369003918
$ws.Range("B21").Value = 12.73128171024789
$ws.Range("C21").Value = 7.926050294395166
$ws.Range("D21").Value = 15.38152165052727
$ws.Range("E21").Value = 16.80813513021495
$ws.Range("G21").Value = 3.669856850090532
$ws.Range("I21").Value = 24.38539706006884
$ws.Range("J21").Value = 9.379954644375962
$ws.Range("K21").Value = 12.88679524453427
$ws.Range("N21").Value = 19.62060182361028
$ws.Range("O21").Value = 27.70627070930373
$ws.Range("B22").Value = 12.94083178555965
$ws.Range("C22").Value = 8.065793856573421
$ws.Range("D22").Value = 15.45704732979528
$ws.Range("E22").Value = 16.88565436957765
$ws.Range("G22").Value = 3.668368793254616
$ws.Range("I22").Value = 24.33999675685267
$ws.Range("J22").Value = 9.381865020658047
$ws.Range("K22").Value = 13.03343868292973
$ws.Range("N22").Value = 19.57397150238915
$ws.Range("O22").Value = 27.68214436948628
$ws.Range("B23").Value = 12.82926486085379
$ws.Range("C23").Value = 7.991413192319671
$ws.Range("D23").Value = 15.41652179906928
$ws.Range("E23").Value = 16.84403641536417
$ws.Range("G23").Value = 3.669157712164711
$ws.Range("I23").Value = 24.36387051910532
$ws.Range("J23").Value = 9.380733774568398
$ws.Range("K23").Value = 12.95521643237822
$ws.Range("N23").Value = 19.59872161688918
$ws.Range("O23").Value = 27.69456900530029
$ws.Range("B24").Value = 12.4011032957496
$ws.Range("C24").Value = 7.705509616328958
$ws.Range("D24").Value = 15.26788423372801
$ws.Range("E24").Value = 16.6918854748892
$ws.Range("G24").Value = 3.672262005708361
$ws.Range("I24").Value = 24.46212741472639
$ws.Range("J24").Value = 9.378906199952784
$ws.Range("K24").Value = 12.65825080007511
$ws.Range("N24").Value = 19.69548340461785
$ws.Range("O24").Value = 27.75154686796754
$ws.Range("B25").Value = 11.92992063816428
$ws.Range("C25").Value = 7.389922585994587
$ws.Range("D25").Value = 15.11852264377608
$ws.Range("E25").Value = 16.5400621898729
$ws.Range("G25").Value = 3.675860946329186
$ws.Range("I25").Value = 24.58480048386985
$ws.Range("J25").Value = 9.382207461703246
$ws.Range("K25").Value = 12.33808012093405
$ws.Range("N25").Value = 19.80636263790058
$ws.Range("O25").Value = 27.83413527690107
